$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-CellText 'D2' '31.604.16'
Set-CellText 'E2' '  +4.16%  '

Set-CellText 'D3' '1.992.20'
Set-CellText 'E3' '  +6.26%  '

Set-CellText 'D4' '0.9978'
Set-CellText 'E4' '  -0.54%  '

Set-CellText 'D5' '0.7890'
Set-CellText 'E5' '  +67.88%  '

Set-CellText 'D6' '253.04'
Set-CellText 'E6' '  +4.00%  '

Set-CellText 'D7' '0.9965'

Set-CellText 'D8' '0.3436'
Set-CellText 'E8' '  +19.93%  '

Set-CellText 'D9' '25.78'
Set-CellText 'E9' '  +17.65%  '

Set-CellText 'D10' '0.06975'
Set-CellText 'E10' '  +8.83%  '

Set-CellText 'D11' '0.8449'
Set-CellText 'E11' '  +17.53%  '

Set-CellText 'D12' '0.08139'
Set-CellText 'E12' '  +4.56%  '

Set-CellText 'D13' '102.62'
Set-CellText 'E13' '  +8.05%  '

Set-CellText 'D14' '1.983.29'
Set-CellText 'E14' '  +4.89%  '

Set-CellText 'D15' '5.514'
Set-CellText 'E15' '  +7.28%  '

Set-CellText 'D16' '277.01'
Set-CellText 'E16' '  -0.31%  '

Set-CellText 'D17' '31.549.53'
Set-CellText 'E17' '  +3.89%  '

Set-CellText 'D18' '14.04'
Set-CellText 'E18' '  +8.49%  '

Set-CellText 'D19' '0.000007905'
Set-CellText 'E19' '  +6.88%  '

Set-CellText 'D20' '2.242.81'
Set-CellText 'E20' '  +5.19%  '

Set-CellText 'D21' '5.689'
Set-CellText 'E21' '  +9.06%  '

Set-CellText 'D22' '0.9957'
Set-CellText 'E22' '  -0.66%  '

Set-CellText 'D23' '0.9984'
Set-CellText 'E23' '  -0.57%  '

Set-CellText 'D24' '6.905'
Set-CellText 'E24' '  +10.79%  '

Set-CellText 'D25' '9.689'
Set-CellText 'E25' '  +7.55%  '

Set-CellText 'D26' '0.1566'
Set-CellText 'E26' '  +62.79%  '

Set-CellText 'D27' '166.56'
Set-CellText 'E27' '  +1.76%  '

Set-CellText 'D28' '19.77'
Set-CellText 'E28' '  +5.94%  '

Set-CellText 'D29' '2.241'
Set-CellText 'E29' '  +19.48%  '

Set-CellText 'E30' '  +6.68%  '

Set-CellText 'D31' '1.353'
Set-CellText 'E31' '  -0.09%  '

Set-CellText 'E32' '  +8.40%  '

Set-CellText 'D33' '4.328'
Set-CellText 'E33' '  +5.57%  '

Set-CellText 'D34' '0.05217'
Set-CellText 'E34' '  +8.41%  '

Set-CellText 'D35' '1.221'
Set-CellText 'E35' '  +9.44%  '

Set-CellText 'D36' '0.7465'
Set-CellText 'E36' '  +8.97%  '

Set-CellText 'E37' '  +3.01%  '

Set-CellText 'D38' '0.9947'
Set-CellText 'E38' '  -0.65%  '

Set-CellText 'D39' '0.01991'
Set-CellText 'E39' '  +6.52%  '

Set-CellText 'D40' '2.911'
Set-CellText 'E40' '  +3.54%  '

Set-CellText 'D41' '6.600'
Set-CellText 'E41' '  +6.04%  '

Set-CellText 'D42' '78.53'
Set-CellText 'E42' '  +6.04%  '

Set-CellText 'D43' '0.4663'
Set-CellText 'E43' '  +10.59%  '

Set-CellText 'D44' '2.086'
Set-CellText 'E44' '  +7.71%  '

Set-CellText 'B45' 'TrustWalletToken'
Set-CellText 'C45' 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-CellText 'D45' '0.8542'
Set-CellText 'E45' '  +3.65%  '

Set-CellText 'B46' 'Quant'
Set-CellText 'C46' 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-CellText 'D46' '105.57'
Set-CellText 'E46' '  +4.77%  '

Set-CellText 'D47' '0.9971'

Set-CellText 'D48' '9.988'
Set-CellText 'E48' '  +4.29%  '

Set-CellText 'D49' '7.527'
Set-CellText 'E49' '  +9.28%  '

Set-CellText 'B50' 'Elrond'
Set-CellText 'C50' 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-CellText 'D50' '36.58'
Set-CellText 'E50' '  +4.59%  '

Set-CellText 'B51' 'Decentraland'
Set-CellText 'C51' 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-CellText 'D51' '0.4273'
Set-CellText 'E51' '  +9.83%  '

Write-Host "Applied cryptos update"
